# Initial commit for Ripon model
# Updates row 2 of Sheet1 from the Newmarket/Sea Battle record to the
# Newbury/Dalamoi record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 44976
$ws.Range("C2").Value = "NEWBURY"
$ws.Range("D2").Value = "GTF"
$ws.Range("E2").Value = "Handicap Chase"
$ws.Range("F2").Value = "2m6½f"
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = "13:30:00"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Dalamoi"
$ws.Range("K2").Value = 128
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = "11-3"
$ws.Range("O2").Value = "Alan Johns"
$ws.Range("P2").Value = "Tim Vaughan"

$ws.Range("T2").Value = ""
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = ""
$ws.Range("W2").Value = ""
$ws.Range("X2").Value = ""
$ws.Range("Y2").Value = 2.88
$ws.Range("Z2").Value = 1

$ws.Range("AB2").Value = 1

$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = "-"
$ws.Range("AG2").Value = "-"
$ws.Range("AH2").Value = "-"
$ws.Range("AI2").Value = ""
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 9
$ws.Range("AL2").Value = 1
$ws.Range("AM2").Value = 42.54
$ws.Range("AN2").Value = 47.98
$ws.Range("AO2").Value = 167
$ws.Range("AP2").Value = 2
$ws.Range("AQ2").Value = 10
$ws.Range("AR2").Value = 22
$ws.Range("AS2").Value = 0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = 0
$ws.Range("AV2").Value = 0
$ws.Range("AW2").Value = 7
$ws.Range("AX2").Value = "NO"
$ws.Range("AY2").Value = "NO"
$ws.Range("AZ2").Value = "YES"
